# Apply updated cryptocurrency market data to the worksheet.
# Values in column D that are plain decimal numbers are written with a
# leading apostrophe so Excel stores them as text (matching the source data's
# inlineStr cell type) instead of converting them to numbers; the cell Style
# is then reset to "Normal" so no stray number-format/quote-prefix style is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.268.19"
$ws.Range("E2").Value = "  +1.68%  "

# Row 3
$ws.Range("D3").Value = "3.809.44"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'701.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.67%  "

# Row 6
$ws.Range("D6").Value = "'173.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7
$ws.Range("D7").Value = "3.807.78"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").Value = "'0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.24%  "

# Row 11
$ws.Range("D11").Value = "'7.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.50%  "

# Row 12
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.67%  "

# Row 14
$ws.Range("D14").Value = "'36.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").Value = "4.450.24"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").Value = "3.798.87"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").Value = "71.255.51"
$ws.Range("E17").Value = "  +1.75%  "

# Row 18
$ws.Range("D18").Value = "'17.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

# Row 19
$ws.Range("D19").Value = "'7.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("E21").Value = "  +6.77%  "

# Row 22
$ws.Range("D22").Value = "'481.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "

# Row 23
$ws.Range("D23").Value = "'0.714"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").Value = "'84.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.96%  "

# Row 25
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "

# Row 26
$ws.Range("D26").Value = "'12.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("D27").Value = "'10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.64%  "

# Row 28
$ws.Range("D28").Value = "'2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
$ws.Range("D29").Value = "3.958.96"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'3.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.73%  "

# Row 31
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").Value = "'7.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.57%  "

# Row 33
$ws.Range("D33").Value = "'2.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.63%  "

# Row 34
$ws.Range("E34").Value = "  +5.34%  "

# Row 35
$ws.Range("D35").Value = "'29.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'9.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.36%  "

# Row 37
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "'1.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "

# Row 38
$ws.Range("E38").Value = "  +1.56%  "

# Row 39
$ws.Range("D39").Value = "'3.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.17%  "

# Row 40
$ws.Range("D40").Value = "'6.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.93%  "

# Row 41
$ws.Range("D41").Value = "'2.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.72%  "

# Row 42
$ws.Range("D42").Value = "'0.987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "

# Row 43
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("E45").Value = "  +16.24%  "

# Row 46
$ws.Range("D46").Value = "'164.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "

# Row 47
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'44.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.67%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'48.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.62%  "

# Row 49
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.303"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

# Row 50
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'416.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.46%  "

# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
